$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to text formatting before writing, so that
# numeric-looking strings (e.g. "1.00") are preserved verbatim as text
# instead of being auto-converted to numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "51.596.64"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "3.021.75"
$ws.Range("E3").Value = "  +2.87%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "378.72"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").Value = "103.29"
$ws.Range("E6").Value = "  +2.79%  "
$ws.Range("D7").Value = "0.546"
$ws.Range("E7").Value = "  +1.53%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.597"
$ws.Range("E9").Value = "  +3.32%  "
$ws.Range("D10").Value = "36.85"
$ws.Range("E10").Value = "  +2.67%  "
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").Value = "0.0861"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").Value = "3.500.61"
$ws.Range("E13").Value = "  +3.10%  "
$ws.Range("D14").Value = "18.55"
$ws.Range("E14").Value = "  +2.02%  "
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("D16").Value = "3.021.40"
$ws.Range("E16").Value = "  +3.12%  "
$ws.Range("D17").Value = "0.983"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").Value = "10.54"
$ws.Range("D19").Value = "51.623.87"
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("D20").Value = "3.04"
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("D21").Value = "12.53"
$ws.Range("E21").Value = "  +1.34%  "
$ws.Range("E22").Value = "  +1.86%  "
$ws.Range("D23").Value = "70.03"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("D24").Value = "267.98"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("E25").Value = "  -2.46%  "
$ws.Range("D26").Value = "8.22"
$ws.Range("E26").Value = "  +3.86%  "
$ws.Range("D27").Value = "7.54"
$ws.Range("E27").Value = "  +6.70%  "
$ws.Range("E28").Value = "  +6.52%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "26.25"
$ws.Range("E30").Value = "  +2.84%  "
$ws.Range("E31").Value = "  +1.19%  "
$ws.Range("D32").Value = "10.33"
$ws.Range("E32").Value = "  +3.28%  "
$ws.Range("D33").Value = "34.35"
$ws.Range("E33").Value = "  +2.75%  "
$ws.Range("D34").Value = "50.51"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").Value = "2.06"
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("D36").Value = "0.0454"
$ws.Range("E36").Value = "  +5.43%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").Value = "3.30"
$ws.Range("E38").Value = "  +6.82%  "
$ws.Range("D39").Value = "17.41"
$ws.Range("E39").Value = "  +5.45%  "
$ws.Range("D40").Value = "0.287"
$ws.Range("E40").Value = "  +10.98%  "
$ws.Range("E41").Value = "  +3.25%  "
$ws.Range("D42").Value = "2.58"
$ws.Range("E42").Value = "  +5.03%  "
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("D44").Value = "127.30"
$ws.Range("E44").Value = "  +6.01%  "
$ws.Range("D45").Value = "3.73"
$ws.Range("E45").Value = "  +9.18%  "
$ws.Range("D46").Value = "21.71"
$ws.Range("E46").Value = "  +3.05%  "
$ws.Range("D47").Value = "2.07"
$ws.Range("E47").Value = "  +2.85%  "
$ws.Range("D48").Value = "2.36"
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("D49").Value = "2.034.62"
$ws.Range("E49").Value = "  +1.39%  "
$ws.Range("D50").Value = "3.321.88"
$ws.Range("E50").Value = "  +3.04%  "
$ws.Range("E51").Value = "  +2.40%  "

# Restore the original (default) cell style now that the text values are set,
# so the cell formatting matches the source workbook.
$priceRange.Style = "Normal"
